$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original text formatting
# (values like "69.577.96" or "575.97" must stay text, not be parsed as numbers)
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.577.96'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.507.55'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.97'
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.01'
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.516'
$ws.Range('E8').Value = '  -0.55%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.506.05'
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.162'
$ws.Range('E10').Value = '  +1.92%  '
$ws.Range('E11').Value = '  -0.56%  '
$ws.Range('E12').Value = '  +6.29%  '
$ws.Range('E13').Value = '  +1.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.967.63'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000178'
$ws.Range('E15').Value = '  +1.40%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.337.47'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.86'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.509.16'
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.31'
$ws.Range('E19').Value = '  -1.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.58'
$ws.Range('E20').Value = '  -2.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '351.74'
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.92'
$ws.Range('E22').Value = '  -0.72%  '
$ws.Range('E23').Value = '  -0.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.73'
$ws.Range('E25').Value = '  +2.27%  '
$ws.Range('E26').Value = '  -1.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.86'
$ws.Range('E27').Value = '  -1.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.665.99'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.01'
$ws.Range('E29').Value = '  +0.75%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0893'
$ws.Range('E30').Value = '  -1.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.87'
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '461.20'
$ws.Range('E32').Value = '  -3.85%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.23'
$ws.Range('E33').Value = '  -5.31%  '
$ws.Range('E34').Value = '  -0.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '159.93'
$ws.Range('E36').Value = '  +4.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.116'
$ws.Range('E37').Value = '  +1.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.07'
$ws.Range('E38').Value = '  +1.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.48'
$ws.Range('E39').Value = '  -0.45%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.319'
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('E42').Value = '  -1.73%  '
$ws.Range('E43').Value = '  -1.52%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.21'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.23'
$ws.Range('E45').Value = '  -4.84%  '
$ws.Range('E46').Value = '  -6.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '142.26'
$ws.Range('E47').Value = '  -1.10%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.48'
$ws.Range('E48').Value = '  -1.96%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.521'
$ws.Range('E49').Value = '  -1.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0735'
$ws.Range('E50').Value = '  +0.66%  '
$ws.Range('E51').Value = '  +2.82%  '
